$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1872909698996655
$ws.Range("C2").Value = 0.5250836120401338
$ws.Range("J2").Value = 0.02341137123745819
$ws.Range("P2").Value = 0.1337792642140468
$ws.Range("S2").Value = 0.1304347826086956
$ws.Range("B3").Value = 0.01666666666666667
$ws.Range("C3").Value = 0.03333333333333333
$ws.Range("J3").Value = 0.07777777777777778
$ws.Range("P3").Value = 0.7111111111111111
$ws.Range("S3").Value = 0.1611111111111111
$ws.Range("J4").Value = 0.1538461538461539
$ws.Range("P4").Value = 0.7307692307692307
$ws.Range("S4").Value = 0.1153846153846154
$ws.Range("B6").Value = 0.06521739130434782
$ws.Range("D6").Value = 0.008695652173913044
$ws.Range("F6").Value = 0.06521739130434782
$ws.Range("J6").Value = 0.2956521739130435
$ws.Range("O6").Value = 0.0391304347826087
$ws.Range("Q6").Value = 0.2521739130434782
$ws.Range("R6").Value = 0.03043478260869565
$ws.Range("S6").Value = 0.2434782608695652
$ws.Range("B7").Value = 0.1397379912663755
$ws.Range("D7").Value = 0.008733624454148471
$ws.Range("F7").Value = 0.06550218340611354
$ws.Range("J7").Value = 0.1746724890829694
$ws.Range("O7").Value = 0.004366812227074236
$ws.Range("Q7").Value = 0.1790393013100437
$ws.Range("R7").Value = 0.0611353711790393
$ws.Range("S7").Value = 0.3668122270742358
$ws.Range("B8").Value = 0.0963855421686747
$ws.Range("D8").Value = 0.02610441767068273
$ws.Range("F8").Value = 0.05622489959839357
$ws.Range("J8").Value = 0.106425702811245
$ws.Range("O8").Value = 0.01606425702811245
$ws.Range("Q8").Value = 0.2650602409638554
$ws.Range("R8").Value = 0.06827309236947791
$ws.Range("S8").Value = 0.3654618473895582
$ws.Range("B9").Value = 0.07096774193548387
$ws.Range("D9").Value = 0.01290322580645161
$ws.Range("F9").Value = 0.09677419354838709
$ws.Range("J9").Value = 0.1870967741935484
$ws.Range("O9").Value = 0.03870967741935484
$ws.Range("Q9").Value = 0.1419354838709677
$ws.Range("R9").Value = 0.05806451612903226
$ws.Range("S9").Value = 0.3935483870967742
$ws.Range("B10").Value = 0.09527326440177253
$ws.Range("D10").Value = 0.02289512555391433
$ws.Range("E10").Value = 0.001477104874446086
$ws.Range("F10").Value = 0.06425406203840472
$ws.Range("J10").Value = 0.1189069423929099
$ws.Range("O10").Value = 0.02880354505169867
$ws.Range("Q10").Value = 0.2548005908419498
$ws.Range("R10").Value = 0.06646971935007386
$ws.Range("S10").Value = 0.3471196454948302
$ws.Range("F11").Value = 0.003012048192771084
$ws.Range("G11").Value = 0.1234939759036145
$ws.Range("J11").Value = 0.0963855421686747
$ws.Range("K11").Value = 0.1807228915662651
$ws.Range("L11").Value = 0.5903614457831325
$ws.Range("S11").Value = 0.006024096385542169
$ws.Range("G12").Value = 0.7598039215686274
$ws.Range("J12").Value = 0.1813725490196078
$ws.Range("K12").Value = 0.009803921568627451
$ws.Range("L12").Value = 0.02941176470588235
$ws.Range("S12").Value = 0.0196078431372549
$ws.Range("G13").Value = 0.7346938775510204
$ws.Range("J13").Value = 0.2448979591836735
$ws.Range("S13").Value = 0.02040816326530612
$ws.Range("F15").Value = 0.01984126984126984
$ws.Range("H15").Value = 0.1746031746031746
$ws.Range("I15").Value = 0.07142857142857142
$ws.Range("J15").Value = 0.3690476190476191
$ws.Range("K15").Value = 0.04365079365079365
$ws.Range("M15").Value = 0.0119047619047619
$ws.Range("O15").Value = 0.02777777777777778
$ws.Range("S15").Value = 0.2817460317460317
$ws.Range("F16").Value = 0.02030456852791878
$ws.Range("H16").Value = 0.1776649746192893
$ws.Range("I16").Value = 0.06091370558375635
$ws.Range("J16").Value = 0.3807106598984771
$ws.Range("K16").Value = 0.1319796954314721
$ws.Range("M16").Value = 0.02030456852791878
$ws.Range("O16").Value = 0.04568527918781726
$ws.Range("S16").Value = 0.1624365482233502
$ws.Range("F17").Value = 0.02010050251256281
$ws.Range("H17").Value = 0.1876046901172529
$ws.Range("I17").Value = 0.05695142378559464
$ws.Range("J17").Value = 0.407035175879397
$ws.Range("K17").Value = 0.09547738693467336
$ws.Range("M17").Value = 0.02177554438860971
$ws.Range("N17").Value = 0.001675041876046901
$ws.Range("O17").Value = 0.07872696817420435
$ws.Range("S17").Value = 0.1306532663316583
$ws.Range("F18").Value = 0.0125
$ws.Range("H18").Value = 0.20625
$ws.Range("I18").Value = 0.0625
$ws.Range("J18").Value = 0.38125
$ws.Range("K18").Value = 0.11875
$ws.Range("M18").Value = 0.00625
$ws.Range("N18").Value = 0.00625
$ws.Range("O18").Value = 0.08125
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.01916932907348243
$ws.Range("H19").Value = 0.2196485623003195
$ws.Range("I19").Value = 0.06469648562300319
$ws.Range("J19").Value = 0.3578274760383386
$ws.Range("K19").Value = 0.1293929712460064
$ws.Range("M19").Value = 0.0231629392971246
$ws.Range("N19").Value = 0.002396166134185303
$ws.Range("O19").Value = 0.07108626198083066
$ws.Range("S19").Value = 0.1126198083067093
